$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new date headers inserted before existing ones (columns shift right by 9)
# Step 1: write the 9 new (most recent) week labels into B1:J1, right-to-left so that
#         new shared-string entries are created in ascending chronological order
#         (Jun_16..Sep_08 -> indices 52..60), matching the target workbook.
$ws.Range("J1").Value2 = "Jun_16"
$ws.Range("I1").Value2 = "Jun_24"
$ws.Range("H1").Value2 = "Jun_30"
$ws.Range("G1").Value2 = "Jul_07"
$ws.Range("F1").Value2 = "Jul_17"
$ws.Range("E1").Value2 = "Jul_23"
$ws.Range("D1").Value2 = "Aug_04"
$ws.Range("C1").Value2 = "Aug_25"
$ws.Range("B1").Value2 = "Sep_08"

# Step 2: shift the previously-existing week labels (old B1:S1) right into K1:AB1
$ws.Range("K1").Value2 = "Jun_09"
$ws.Range("L1").Value2 = "Jun_03"
$ws.Range("M1").Value2 = "May_27"
$ws.Range("N1").Value2 = "May_23"
$ws.Range("O1").Value2 = "May_19"
$ws.Range("P1").Value2 = "May_15"
$ws.Range("Q1").Value2 = "May_12"
$ws.Range("R1").Value2 = "May_05"
$ws.Range("S1").Value2 = "Apr_28"
$ws.Range("T1").Value2 = "Apr_24"
$ws.Range("U1").Value2 = "Apr_21"
$ws.Range("V1").Value2 = "Apr_17"
$ws.Range("W1").Value2 = "Apr_11"
$ws.Range("X1").Value2 = "Apr_07"
$ws.Range("Y1").Value2 = "Apr_04"
$ws.Range("Z1").Value2 = "Mar_31"
$ws.Range("AA1").Value2 = "Mar_27"
$ws.Range("AB1").Value2 = "Mar_24"

# Data rows: append 9 new "UN" placeholder cells after each row's existing last column
# (existing cells are left untouched/in place, matching the source workbook's pattern).
foreach ($r in 2..29) {
    foreach ($colNum in 20..28) {
        $ws.Cells.Item($r, $colNum).Value2 = "UN"
    }
}
foreach ($r in 30..31) {
    foreach ($colNum in 17..25) {
        $ws.Cells.Item($r, $colNum).Value2 = "UN"
    }
}
foreach ($r in 32..33) {
    foreach ($colNum in 8..16) {
        $ws.Cells.Item($r, $colNum).Value2 = "UN"
    }
}
